# feat: add single and multi corrector
# Shift columns left (drop old "NO" column A), update values, drop duplicate
# footer cell, and adjust the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the numeric-looking result values (e.g. "7.08") as plain text,
# matching the source data which stores everything as inline strings.
$ws.Range("A1:D4").NumberFormat = "@"
$ws.Range("A5").NumberFormat = "@"

# Row 1 - headers
$ws.Range("A1").Value = "项目"
$ws.Range("B1").Value = "结果"
$ws.Range("C1").Value = "参考值"
$ws.Range("D1").Value = "单位"

# Row 2
$ws.Range("A2").Value = "促甲状腺激素"
$ws.Range("B2").Value = "7.08"
$ws.Range("C2").Value = "0.27-4.20"
$ws.Range("D2").Value = "IU/mL"

# Row 3
$ws.Range("A3").Value = "游离甲状腺素"
$ws.Range("B3").Value = "18.15"
$ws.Range("C3").Value = "12.0-22.0"
$ws.Range("D3").Value = "mol/L"

# Row 4
$ws.Range("A4").Value = "游离三碘甲状腺原氨酸"
$ws.Range("B4").Value = "5.11"
$ws.Range("C4").Value = "3.10-6.80"
$ws.Range("D4").Value = "mol/L"

# Row 5 - footer note, only column A now (drop the duplicate B5 cell
# entirely, it should no longer exist)
$ws.Range("A5").Value = "如对检验结果咨询，请在48小时内与检验科联系"
$ws.Range("B5").ClearContents()

# Clear the old column E which is no longer part of the table
$ws.Range("E1:E5").ClearContents()
